$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("Q2").Value = 2.35
$ws.Range("R2").Value = 1.57

# Row 6 updates
$ws.Range("G6").Value = 1.32
$ws.Range("H6").Value = 4.8
$ws.Range("K6").Value = 2.55
$ws.Range("L6").Value = 6.7
$ws.Range("N6").Value = 9.25
$ws.Range("O6").Value = 1.16
$ws.Range("P6").Value = 4.5
$ws.Range("Q6").Value = 1.5
$ws.Range("R6").Value = 2.4
$ws.Range("S6").Value = 1.28
$ws.Range("T6").Value = 3.35
$ws.Range("U6").Value = 1.8
$ws.Range("V6").Value = 1.91
$ws.Range("W6").Value = 8.75
$ws.Range("X6").Value = 7.3
$ws.Range("Z6").Value = 9
$ws.Range("AB6").Value = 23
$ws.Range("AC6").Value = 9.25
$ws.Range("AD6").Value = 10
$ws.Range("AE6").Value = 19
$ws.Range("AF6").Value = 75
$ws.Range("AG6").Value = 23
$ws.Range("AH6").Value = 55
$ws.Range("AL6").Value = 65
$ws.Range("AM6").Value = 500
$ws.Range("AN6").Value = 3.3
$ws.Range("AO6").Value = 5.9
$ws.Range("AP6").Value = 14.5
$ws.Range("AQ6").Value = 15
$ws.Range("AR6").Value = 37
$ws.Range("AS6").Value = 175
$ws.Range("AT6").Value = 3.35
$ws.Range("AV6").Value = 65
$ws.Range("AX6").Value = 45
$ws.Range("AY6").Value = 40
$ws.Range("AZ6").Value = 300
$ws.Range("BA6").Value = 300
$ws.Range("BB6").Value = 500
